$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.817.25"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "3.298.69"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").Value = "3.293.15"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.570"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.175"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.570"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.56%  "

$ws.Range("D14").Value = "3.834.76"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "611.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.56%  "

$ws.Range("D17").Value = "65.841.87"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "

$ws.Range("D20").Value = "3.309.43"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.888"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "550.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.54%  "

$ws.Range("D35").Value = "3.799.98"
$ws.Range("E35").Value = "  -0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.128"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.84%  "

$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.98%  "

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0676"
$ws.Range("E44").Value = "  -8.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0405"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.18%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.126"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.72%  "

